$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E3").Value = 16.297
$ws.Range("E4").Value = 16.499
$ws.Range("C7").Value = -13.566
$ws.Range("B8").Value = 6.779999999999999
$ws.Range("B10").Value = 6.537999999999999
$ws.Range("E11").Value = 16.529
$ws.Range("B12").Value = 5.165999999999999
$ws.Range("C14").Value = -13.205
$ws.Range("E14").Value = 16.752
$ws.Range("C15").Value = -13.636
$ws.Range("B18").Value = 5.762
$ws.Range("C18").Value = -12.379
$ws.Range("E18").Value = 17.157
$ws.Range("E19").Value = 16.452
$ws.Range("C20").Value = -11.93
$ws.Range("E21").Value = 16.547
$ws.Range("B25").Value = 6.145
$ws.Range("E27").Value = 16.584
$ws.Range("C29").Value = -12.334
$ws.Range("C30").Value = -13.177
$ws.Range("C31").Value = -13.41
$ws.Range("E31").Value = 16.151
$ws.Range("C35").Value = -11.986
$ws.Range("B37").Value = 7.996999999999998
$ws.Range("E38").Value = 16.644
$ws.Range("C40").Value = -12.782
$ws.Range("E42").Value = 16.456
$ws.Range("C44").Value = -12.406
$ws.Range("E44").Value = 16.749
$ws.Range("E47").Value = 16.234
$ws.Range("C50").Value = -13.371
$ws.Range("C54").Value = -12.496
$ws.Range("B55").Value = 5.257
$ws.Range("E56").Value = 16.206
$ws.Range("E58").Value = 16.541
$ws.Range("E65").Value = 16.987
$ws.Range("B68").Value = 5.431
$ws.Range("C68").Value = -11.177
$ws.Range("E73").Value = 16.577
$ws.Range("C76").Value = -13.366
$ws.Range("B77").Value = 5.737
$ws.Range("B78").Value = 7.544999999999999
$ws.Range("B79").Value = 5.353
$ws.Range("B80").Value = 8.087999999999999
$ws.Range("B81").Value = 6.154000000000001
$ws.Range("B82").Value = 6.009
$ws.Range("B84").Value = 6.029
$ws.Range("C87").Value = -13.221
$ws.Range("C88").Value = -12.774
$ws.Range("E90").Value = 16.398
$ws.Range("C92").Value = -11.248
$ws.Range("E92").Value = 18.025
$ws.Range("E94").Value = 18.031
$ws.Range("E95").Value = 17.239
$ws.Range("C96").Value = -12.85
$ws.Range("C98").Value = -12.888
$ws.Range("B101").Value = 8.928999999999998
$ws.Range("C101").Value = -12.624
$ws.Range("E101").Value = 16.539
$ws.Range("B102").Value = 7.264
$ws.Range("C102").Value = -12.415
